$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constant")

# Insert a new row above current row 3, shifting rows 3-53 down to 4-54.
# The new row picks up the formatting of the row above (row 2), which matches
# the desired style (fill style index 2) for the new "bg_mixing" row.
$ws.Rows("3").Insert()

# Fill in the newly inserted row 3 (bg_mixing)
$ws.Range("A3").Value = "bg_mixing"
$ws.Range("B3").Value = 0.05
$ws.Range("C3").Value = "uniform"
$ws.Range("D3").Value = 0.01
$ws.Range("E3").Value = 0.1
$ws.Range("G3").Value = "Background age-agnostic mixing level"
$ws.Range("H3").Clear()

# Update row 4 (previously row 3: child_socialising) to a_spread
$ws.Range("A4").Value = "a_spread"
$ws.Range("B4").Value = 10
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 15
$ws.Range("G4").Value = "Spread of assortative mixing pattern (smaller value means more assortativity)"

# Update row 5 (previously row 4: elderly_socialising) to pc_strength
$ws.Range("A5").Value = "pc_strength"
$ws.Range("B5").Value = 1.5
$ws.Range("D5").Value = 0.5
$ws.Range("E5").Value = 5
$ws.Range("G5").Value = "Strength of parent-children mixing pattern"

$ws.Range("G9").Select()
